# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" for the 254acf12-... file across
# the Overview sheet and the two locale sheets to reflect the newly
# generated handoff report timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-29-18 10:29:08"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-18 10:29:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-18 10:29:08"
